$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Phase 1: replace each original value with a unique fixed-width placeholder
# token (no token is a prefix of another) to avoid any collisions where a
# new value equals another cells old value, or vice versa.
Replace-Text "2024-12-19 Thursday" "@@TOK_000@@"
Replace-Text "53-35=18" "@@TOK_001@@"
Replace-Text "49+46=95" "@@TOK_002@@"
Replace-Text "70-64=6" "@@TOK_003@@"
Replace-Text "65+19=84" "@@TOK_004@@"
Replace-Text "97-48=49" "@@TOK_005@@"
Replace-Text "91-55=36" "@@TOK_006@@"
Replace-Text "98-59=39" "@@TOK_007@@"
Replace-Text "91-52=39" "@@TOK_008@@"
Replace-Text "29+16=45" "@@TOK_009@@"
Replace-Text "16+19=35" "@@TOK_010@@"
Replace-Text "18+24=42" "@@TOK_011@@"
Replace-Text "4+9=13" "@@TOK_012@@"
Replace-Text "70-31=39" "@@TOK_013@@"
Replace-Text "56+8=64" "@@TOK_014@@"
Replace-Text "17+57=74" "@@TOK_015@@"
Replace-Text "45-18=27" "@@TOK_016@@"
Replace-Text "49+35=84" "@@TOK_017@@"
Replace-Text "61-49=12" "@@TOK_018@@"
Replace-Text "42-36=6" "@@TOK_019@@"
Replace-Text "44-18=26" "@@TOK_020@@"
Replace-Text "57+17=74" "@@TOK_021@@"
Replace-Text "28+5=33" "@@TOK_022@@"
Replace-Text "41-15=26" "@@TOK_023@@"
Replace-Text "51-49=2" "@@TOK_024@@"
Replace-Text "58+23=81" "@@TOK_025@@"
Replace-Text "19+27=46" "@@TOK_026@@"
Replace-Text "91-4=87" "@@TOK_027@@"
Replace-Text "74-56=18" "@@TOK_028@@"
Replace-Text "35+49=84" "@@TOK_029@@"
Replace-Text "59+19=78" "@@TOK_030@@"
Replace-Text "35-7=28" "@@TOK_031@@"
Replace-Text "56-29=27" "@@TOK_032@@"
Replace-Text "82-44=38" "@@TOK_033@@"
Replace-Text "91-23=68" "@@TOK_034@@"
Replace-Text "36+58=94" "@@TOK_035@@"
Replace-Text "26+57=83" "@@TOK_036@@"
Replace-Text "70-4=66" "@@TOK_037@@"
Replace-Text "44-29=15" "@@TOK_038@@"
Replace-Text "42+49=91" "@@TOK_039@@"
Replace-Text "73-14=59" "@@TOK_040@@"
Replace-Text "29+62=91" "@@TOK_041@@"
Replace-Text "73-44=29" "@@TOK_042@@"
Replace-Text "29+9=38" "@@TOK_043@@"
Replace-Text "32+19=51" "@@TOK_044@@"
Replace-Text "7+85=92" "@@TOK_045@@"
Replace-Text "3+48=51" "@@TOK_046@@"
Replace-Text "14+57=71" "@@TOK_047@@"
Replace-Text "33+18=51" "@@TOK_048@@"
Replace-Text "72-29=43" "@@TOK_049@@"
Replace-Text "55-8=47" "@@TOK_050@@"
Replace-Text "34+27=61" "@@TOK_051@@"
Replace-Text "67-19=48" "@@TOK_052@@"
Replace-Text "74-67=7" "@@TOK_053@@"
Replace-Text "92-79=13" "@@TOK_054@@"
Replace-Text "25-9=16" "@@TOK_055@@"
Replace-Text "84-48=36" "@@TOK_056@@"
Replace-Text "9+26=35" "@@TOK_057@@"
Replace-Text "68-59=9" "@@TOK_058@@"
Replace-Text "9+34=43" "@@TOK_059@@"
Replace-Text "76-39=37" "@@TOK_060@@"
Replace-Text "24-16=8" "@@TOK_061@@"
Replace-Text "44+28=72" "@@TOK_062@@"
Replace-Text "82-4=78" "@@TOK_063@@"
Replace-Text "8+54=62" "@@TOK_064@@"
Replace-Text "41-22=19" "@@TOK_065@@"
Replace-Text "74-57=17" "@@TOK_066@@"
Replace-Text "8+33=41" "@@TOK_067@@"
Replace-Text "66+29=95" "@@TOK_068@@"
Replace-Text "23-6=17" "@@TOK_069@@"
Replace-Text "41-23=18" "@@TOK_070@@"
Replace-Text "35-18=17" "@@TOK_071@@"
Replace-Text "44-16=28" "@@TOK_072@@"
Replace-Text "81-29=52" "@@TOK_073@@"
Replace-Text "26+19=45" "@@TOK_074@@"
Replace-Text "70-47=23" "@@TOK_075@@"
Replace-Text "17+14=31" "@@TOK_076@@"
Replace-Text "25+69=94" "@@TOK_077@@"
Replace-Text "38+38=76" "@@TOK_078@@"
Replace-Text "68+25=93" "@@TOK_079@@"
Replace-Text "54-26=28" "@@TOK_080@@"
Replace-Text "81-24=57" "@@TOK_081@@"
Replace-Text "97-49=48" "@@TOK_082@@"
Replace-Text "70-65=5" "@@TOK_083@@"
Replace-Text "60-15=45" "@@TOK_084@@"
Replace-Text "83-39=44" "@@TOK_085@@"
Replace-Text "27+14=41" "@@TOK_086@@"
Replace-Text "37+57=94" "@@TOK_087@@"
Replace-Text "5+9=14" "@@TOK_088@@"
Replace-Text "36+59=95" "@@TOK_089@@"
Replace-Text "84-7=77" "@@TOK_090@@"
Replace-Text "13-5=8" "@@TOK_091@@"
Replace-Text "38+26=64" "@@TOK_092@@"
Replace-Text "19+35=54" "@@TOK_093@@"
Replace-Text "46+17=63" "@@TOK_094@@"
Replace-Text "60-51=9" "@@TOK_095@@"
Replace-Text "17+27=44" "@@TOK_096@@"
Replace-Text "29+68=97" "@@TOK_097@@"
Replace-Text "3+79=82" "@@TOK_098@@"
Replace-Text "26+56=82" "@@TOK_099@@"
Replace-Text "39+43=82" "@@TOK_100@@"

# Phase 2: replace each placeholder token with the final value
Replace-Text "@@TOK_000@@" "2024-12-20 Friday"
Replace-Text "@@TOK_001@@" "48+8=56"
Replace-Text "@@TOK_002@@" "80-55=25"
Replace-Text "@@TOK_003@@" "91-28=63"
Replace-Text "@@TOK_004@@" "27+35=62"
Replace-Text "@@TOK_005@@" "25+27=52"
Replace-Text "@@TOK_006@@" "63-28=35"
Replace-Text "@@TOK_007@@" "73-48=25"
Replace-Text "@@TOK_008@@" "3+59=62"
Replace-Text "@@TOK_009@@" "14+7=21"
Replace-Text "@@TOK_010@@" "36+9=45"
Replace-Text "@@TOK_011@@" "83-49=34"
Replace-Text "@@TOK_012@@" "63-36=27"
Replace-Text "@@TOK_013@@" "16+28=44"
Replace-Text "@@TOK_014@@" "70-45=25"
Replace-Text "@@TOK_015@@" "68+18=86"
Replace-Text "@@TOK_016@@" "7+8=15"
Replace-Text "@@TOK_017@@" "97-88=9"
Replace-Text "@@TOK_018@@" "97-29=68"
Replace-Text "@@TOK_019@@" "20-14=6"
Replace-Text "@@TOK_020@@" "66-18=48"
Replace-Text "@@TOK_021@@" "41-19=22"
Replace-Text "@@TOK_022@@" "18+65=83"
Replace-Text "@@TOK_023@@" "25+26=51"
Replace-Text "@@TOK_024@@" "74-17=57"
Replace-Text "@@TOK_025@@" "60-59=1"
Replace-Text "@@TOK_026@@" "50-43=7"
Replace-Text "@@TOK_027@@" "94-18=76"
Replace-Text "@@TOK_028@@" "31-13=18"
Replace-Text "@@TOK_029@@" "51-48=3"
Replace-Text "@@TOK_030@@" "27+9=36"
Replace-Text "@@TOK_031@@" "73+18=91"
Replace-Text "@@TOK_032@@" "65-29=36"
Replace-Text "@@TOK_033@@" "97-29=68"
Replace-Text "@@TOK_034@@" "26+45=71"
Replace-Text "@@TOK_035@@" "26+16=42"
Replace-Text "@@TOK_036@@" "18+77=95"
Replace-Text "@@TOK_037@@" "82-14=68"
Replace-Text "@@TOK_038@@" "70-6=64"
Replace-Text "@@TOK_039@@" "53+19=72"
Replace-Text "@@TOK_040@@" "37+59=96"
Replace-Text "@@TOK_041@@" "7+44=51"
Replace-Text "@@TOK_042@@" "9+14=23"
Replace-Text "@@TOK_043@@" "24+57=81"
Replace-Text "@@TOK_044@@" "74-18=56"
Replace-Text "@@TOK_045@@" "81-59=22"
Replace-Text "@@TOK_046@@" "19+9=28"
Replace-Text "@@TOK_047@@" "68+9=77"
Replace-Text "@@TOK_048@@" "16+18=34"
Replace-Text "@@TOK_049@@" "63-59=4"
Replace-Text "@@TOK_050@@" "13+38=51"
Replace-Text "@@TOK_051@@" "34-5=29"
Replace-Text "@@TOK_052@@" "65-17=48"
Replace-Text "@@TOK_053@@" "80-38=42"
Replace-Text "@@TOK_054@@" "28+45=73"
Replace-Text "@@TOK_055@@" "86-79=7"
Replace-Text "@@TOK_056@@" "65+16=81"
Replace-Text "@@TOK_057@@" "5+8=13"
Replace-Text "@@TOK_058@@" "59+3=62"
Replace-Text "@@TOK_059@@" "63-47=16"
Replace-Text "@@TOK_060@@" "27+4=31"
Replace-Text "@@TOK_061@@" "60-4=56"
Replace-Text "@@TOK_062@@" "22-15=7"
Replace-Text "@@TOK_063@@" "58+25=83"
Replace-Text "@@TOK_064@@" "73+8=81"
Replace-Text "@@TOK_065@@" "37+25=62"
Replace-Text "@@TOK_066@@" "29+18=47"
Replace-Text "@@TOK_067@@" "69+6=75"
Replace-Text "@@TOK_068@@" "37+55=92"
Replace-Text "@@TOK_069@@" "71-32=39"
Replace-Text "@@TOK_070@@" "93-39=54"
Replace-Text "@@TOK_071@@" "91-85=6"
Replace-Text "@@TOK_072@@" "18+6=24"
Replace-Text "@@TOK_073@@" "62-28=34"
Replace-Text "@@TOK_074@@" "29+3=32"
Replace-Text "@@TOK_075@@" "25+6=31"
Replace-Text "@@TOK_076@@" "74-48=26"
Replace-Text "@@TOK_077@@" "28+38=66"
Replace-Text "@@TOK_078@@" "23+19=42"
Replace-Text "@@TOK_079@@" "17+55=72"
Replace-Text "@@TOK_080@@" "86-8=78"
Replace-Text "@@TOK_081@@" "38+45=83"
Replace-Text "@@TOK_082@@" "51-26=25"
Replace-Text "@@TOK_083@@" "97-49=48"
Replace-Text "@@TOK_084@@" "32+9=41"
Replace-Text "@@TOK_085@@" "96-7=89"
Replace-Text "@@TOK_086@@" "86-17=69"
Replace-Text "@@TOK_087@@" "70-34=36"
Replace-Text "@@TOK_088@@" "68+5=73"
Replace-Text "@@TOK_089@@" "75+18=93"
Replace-Text "@@TOK_090@@" "17+64=81"
Replace-Text "@@TOK_091@@" "73-59=14"
Replace-Text "@@TOK_092@@" "96-78=18"
Replace-Text "@@TOK_093@@" "60-59=1"
Replace-Text "@@TOK_094@@" "16+55=71"
Replace-Text "@@TOK_095@@" "8+78=86"
Replace-Text "@@TOK_096@@" "37-9=28"
Replace-Text "@@TOK_097@@" "71-19=52"
Replace-Text "@@TOK_098@@" "8+65=73"
Replace-Text "@@TOK_099@@" "69+6=75"
Replace-Text "@@TOK_100@@" "85-48=37"

Write-Output "done"